$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header: "en_with" -> "en_width"
$ws.Range("B1").Value = "en_width"

# Update the active cell selection shown in the sheet view
$ws.Range("D11").Select()
